$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Rarity" column (C) to hold the new
# "Card Type (Weapons)" data. This shifts the old C/D/E (Rarity/Ability/
# CardDescribe) columns to D/E/F.
$ws.Columns("C:C").Insert()

# New column header.
$ws.Range("C1").Value = "Card Type (Weapons)"

# Existing rows (2-5) just need the new Weapons card-type value filled in.
$ws.Range("C2:C5").Value = "Weapons"

# New weapon cards (rows 6-15).
$weapons = @("AK-47", "Deagle", "Glock-18", "HE", "M4A4", "PSL", "Sawed", "Smoke", "UMP", "UZI")
$row = 6
foreach ($weapon in $weapons) {
    $ws.Range("A$row").Formula = "=ROW()-2"
    $ws.Range("B$row").Value = $weapon
    $ws.Range("C$row").Value = "Weapons"
    $ws.Range("D$row").Value = 0
    $ws.Range("E$row").Value = "None"
    $ws.Range("F$row").Value = "None"
    $row = $row + 1
}

$ws.Range("E6:F15").Select()
